$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ID=1, Shopkeeper_ID=20, Brand=Candyland, Total_Due=-920, Last_Payment_Date=2025-03-23 22:07:30, IsDeleted=0
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = "Candyland"
$ws.Range("D2").Value = -920
$ws.Range("E2").Value = "2025-03-23 22:07:30"
$ws.Range("F2").Value = 0

# Row 3: ID=2, Shopkeeper_ID=16, Brand=Bonapapa, Total_Due=-718, Last_Payment_Date=2025-03-23 22:14:37, IsDeleted=0
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 16
$ws.Range("C3").Value = "Bonapapa"
$ws.Range("D3").Value = -718
$ws.Range("E3").Value = "2025-03-23 22:14:37"
$ws.Range("F3").Value = 0

# Row 4: ID=3, Shopkeeper_ID=17, Brand=Bonapapa, Total_Due=0, Last_Payment_Date=(empty), IsDeleted=0
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 17
$ws.Range("C4").Value = "Bonapapa"
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = 0

# Row 5: ID=4, Shopkeeper_ID=18, Brand=Candyland, Total_Due=0, Last_Payment_Date=(empty), IsDeleted=0
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 18
$ws.Range("C5").Value = "Candyland"
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 0

# Row 6 (new): ID=5, Shopkeeper_ID=19, Brand=Candyland, Total_Due=0, Last_Payment_Date=(empty), IsDeleted=0
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 19
$ws.Range("C6").Value = "Candyland"
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = 0
